$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns: EB = study_name_drug, EC = crossover_periods
$ws.Range("EB1").Value = "study_name_drug"
$ws.Range("EC1").Value = "crossover_periods"

# Copy the header style (bold, centered) from the existing last header cell (EA1)
$ws.Range("EA1").Copy() | Out-Null
$ws.Range("EB1:EC1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Per-row values for study_name_drug (EB) and crossover_periods (EC)
$rowData = @{
    2  = @("Isaacson (2023) - ulotaront", 1)
    3  = @("Isaacson (2023) - ulotaront", 1)
    4  = @("Koblan (2020) - ulotaront", 1)
    5  = @("Koblan (2020) - ulotaront", 1)
    6  = @("NCT04072354 (2019) - ulotaront", 1)
    7  = @("NCT04072354 (2019) - ulotaront", 1)
    8  = @("NCT04072354 (2019) - ulotaront", 1)
    9  = @("NCT04092686 (2019) - ulotaront", 1)
    10 = @("NCT04092686 (2019) - ulotaront", 1)
    11 = @("NCT04092686 (2019) - ulotaront", 1)
    12 = @("NCT04512066 (2020) - ralmitaront", 1)
    13 = @("NCT04512066 (2020) - ralmitaront", 1)
    14 = @("NCT04512066 (2020) - ralmitaront", 1)
    15 = @("NCT04512066 (2020) - ralmitaront", 1)
    16 = @("Perini (2023) - ulotaront", 1)
    17 = @("Perini (2023) - ulotaront", 1)
    18 = @("Perini (2023) - ulotaront", 1)
    19 = @("Tsukada (2023) - ulotaront", 2)
    20 = @("Tsukada (2023) - ulotaront", 2)
    21 = @("Hopkins (2021) - ulotaront", 2)
    22 = @("Hopkins (2021) - ulotaront", 2)
    23 = @("Hopkins (2021) - ulotaront", 2)
    24 = @("Szabo (2023) - ulotaront", 3)
    25 = @("Szabo (2023) - ulotaront", 3)
    26 = @("Szabo (2023) - ulotaront", 3)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 132).Value = $vals[0]   # EB -> column 132
    $ws.Cells.Item($r, 133).Value = $vals[1]   # EC -> column 133
}
